# Adds more Human Detection test-result rows to Sheet1 (rows 10-37),
# matching the "added more test results" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns: A=S.No, B=Source position1, C=Source position2,
#          D=Target position, E=Target Velocity, F=h1, G=h2, H=hres,
#          I=Object detected
$rows = @(
    @(10, 8.0, '0,5,0', '0,5.5,0', '4,8,0', '2,1,1', 0.011109, 0.0176, 0.0006585, 'yes'),
    @(11, 9.0, '0,5,0', '0,5.5,0', '15,8,0', '2,1,0', 0.000007447, 0.001723, 0.00006049, 'no'),
    @(12, 10.0, '0,5,0', '0,5.5,0', '10,8,0', '2,1,0', 0.00396, 0.0003597, 0.0001543, 'yes'),
    @(13, 11.0, '0,5,0', '0,6,0', '10,8,0', '2,1,0', 0.00396, 0.004059, 0.0001446, 'yes'),
    @(14, 12.0, '0,5,0', '0,5.5,0', '6,4,0', '2,1,0', 0.00423, 0.01145, 0.0005714, 'yes'),
    @(15, 13.0, '0,5,0', '0,5.5,0', '6,4,1', '3,5,0', 0.00423, 0.01145, 0.008435, 'yes'),
    @(16, 14.0, '0,5,0', '0,5.5,0', '10,10,2', '1,1,0', 0.002257, 0.003309, 0.00004607, 'no'),
    @(17, 15.0, '0,5,0', '0,5.5,0', '10,10,2', '2,3,0', 0.002257, 0.003309, 0.0003295, 'yes'),
    @(18, 16.0, '0,5,0', '0,5.5,0', '10,10,2', '5,2,0', 0.002257, 0.003309, 0.00001643, 'no'),
    @(19, 17.0, '0,5,0', '0,5.5,0', '10,10,2', '2,4,0', 0.002257, 0.003309, 0.0003364, 'yes'),
    @(20, 18.0, '0,5,0', '0,5.5,0', '10,10,2', '2,8,0', 0.002257, 0.003309, 0.0006598, 'yes'),
    @(21, 19.0, '0,5,0', '0,5.5,0', '5,5,0', '1,0,0', 0.0111092, 0.006309, 0.000597, 'yes'),
    @(22, 20.0, '0,5,0', '0,5.5,0', '5,5,0', '10,0,0', 0.0111092, 0.006309, 0.002154, 'yes'),
    @(23, 21.0, '0,5,0', '0,5.5,0', '8,5,0', '1,0,0', 0.003049, 0.004199, 0.0001462, 'yes'),
    @(24, 22.0, '0,5,0', '0,5.5,0', '10,5,0', '1,0,0', 0.002626, 0.002039, 0.00007488, 'no'),
    @(25, 23.0, '0,5,0', '0,5.5,0', '10,5,0', '2,0,0', 0.002626, 0.002039, 0.0002072, 'yes'),
    @(26, 24.0, '0,5,0', '0,5.5,0', '12,5,0', '2,0,0', 0.001781, 0.001441, 0.0001209, 'yes'),
    @(27, 25.0, '0,5,0', '0,5.5,0', '15,5,0', '2,0,0', 0.001098, 0.0009227, 0.00006269, 'no'),
    @(28, 26.0, '0,5,0', '0,5.5,0', '15,5,0', '6,0,0', 0.001098, 0.0009227, 0.00005937, 'no'),
    @(29, 27.0, '0,5,0', '0,5.5,0', '15,5,0', '6,6,0', 0.001098, 0.0009227, 0.003998, 'yes'),
    @(30, 28.0, '0,5,0', '0,5.5,0', '13,5,0', '2,0,0', 0.001499, 0.001231, 0.00009546, 'no'),
    @(31, 29.0, '0,5,0', '0,5.5,0', '5,10,0', '0,2,0', 0.001969, 0.007955, 0.001079, 'yes'),
    @(32, 30.0, '0,5,0', '0,5.5,0', '5,15,0', '0,2,0', 0.003504, 0.002779, 0.0003913, 'yes'),
    @(33, 31.0, '0,5,0', '0,5.5,0', '5,17,0', '0,2,0', 0.001499, 0.0003272, 0.0001014, 'yes'),
    @(34, 32.0, '0,5,0', '0,5.5,0', '5,20,0', '0,2,0', 0.0003153, 0.0009096, 0.0001129, 'yes'),
    @(35, 33.0, '0,5,0', '0,5.5,0', '5,25,0', '0,2,0', 0.00081, 0.0009609, 0.00005212, 'no'),
    @(36, 34.0, '0,5,0', '0,5.5,0', '5,21,0', '0,2,0', 0.0006349, 0.000214, 0.00005181, 'no'),
    @(37, 35.0, '0,5,0', '0,5.5,0', '5,21,0', '2,8,0', 0.0006349, 0.000214, 0.00001842, 'no')
)

foreach ($row in $rows) {
    $r = $row[0]
    for ($col = 1; $col -le 9; $col++) {
        $value = $row[$col]
        if ($null -ne $value) {
            $ws.Cells.Item($r, $col).Value = $value
        }
    }
}

# Update the view state to match: scrolled so row 14 is at the top,
# with F37 as the active selection (the last-entered cell).
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("F37").Select()
